# repull data, push all data, mean calculation
# Update the "dSF" column (F) values for the rows whose raw data was
# repulled / recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -3
    "F5"  = -4
    "F6"  = 3
    "F7"  = 1
    "F11" = 0
    "F15" = -2
    "F17" = 3
    "F25" = -1
    "F26" = 1
    "F35" = 4
    "F39" = 1
    "F47" = 0
    "F58" = -1
    "F63" = 1
    "F69" = -1
    "F73" = 0
    "F74" = 0
    "F75" = 1
    "F84" = 2
    "F93" = 1
    "F94" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
